$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Select()

$ws.Rows.Item(6).Resize(2).Insert()

$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

Write-Host "Done"
Write-Host $ws.Range("A8").Value
Write-Host $ws.Range("B8").Value
